# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# with freshly scraped values. Values that look like plain numbers
# are prefixed with a leading apostrophe so Excel stores them as
# text (preserving formatting such as trailing zeros) instead of
# auto-converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.782.48"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.294.37"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'97.49"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "'270.02"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'0.609"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "'45.42"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").Value = "'0.0937"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "'7.88"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  +1.82%  "
$ws.Range("D14").Value = "'15.80"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "2.637.88"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "'0.861"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "2.293.49"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "43.783.39"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").Value = "'72.25"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E22").Value = "  +9.14%  "
$ws.Range("D23").Value = "'233.28"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").Value = "'9.09"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("D25").Value = "'2.75"
$ws.Range("E25").Value = "  +8.77%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'11.32"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'38.36"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'176.50"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  +1.08%  "
$ws.Range("D36").Value = "'4.71"
$ws.Range("E36").Value = "  +7.51%  "
$ws.Range("D37").Value = "'0.110"
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("D38").Value = "'0.0352"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").Value = "'3.47"
$ws.Range("E39").Value = "  +3.90%  "
$ws.Range("D40").Value = "'0.237"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "'1.36"
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("D43").Value = "'12.19"
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").Value = "'64.75"
$ws.Range("E44").Value = "  +4.70%  "
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "'8.75"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "'98.61"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").Value = "'0.441"
$ws.Range("E50").Value = "  +7.02%  "
$ws.Range("E51").Value = "  +10.82%  "
